$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new D (Price) / E (Volume 1h) values.
# Only rows/columns present in the diff are updated; everything else stays untouched.
$updates = @(
    @{ Row = 2;  D = "65.734.19";  E = "  +1.71%  " }
    @{ Row = 3;  D = "3.484.75" }
    @{ Row = 5;  D = "580.50";     E = "  +0.35%  " }
    @{ Row = 6;  D = "161.58";     E = "  +3.20%  " }
    @{ Row = 7;  E = "  -0.02%  " }
    @{ Row = 8;  D = "3.487.44";   E = "  +0.59%  " }
    @{ Row = 9;  D = "0.586";      E = "  +5.52%  " }
    @{ Row = 11; D = "0.126";      E = "  +0.39%  " }
    @{ Row = 12; D = "0.444";      E = "  -0.79%  " }
    @{ Row = 13; D = "4.089.62";   E = "  +0.73%  " }
    @{ Row = 14; E = "  -1.55%  " }
    @{ Row = 15; D = "0.0000195";  E = "  -1.22%  " }
    @{ Row = 16; D = "28.63";      E = "  +2.67%  " }
    @{ Row = 17; D = "65.740.27";  E = "  +1.66%  " }
    @{ Row = 18; D = "3.487.04";   E = "  +0.85%  " }
    @{ Row = 19; D = "6.45" }
    @{ Row = 20; D = "14.33";      E = "  -0.50%  " }
    @{ Row = 21; D = "392.47";     E = "  -1.36%  " }
    @{ Row = 22; D = "8.26";       E = "  -3.76%  " }
    @{ Row = 23; D = "0.551";      E = "  +0.45%  " }
    @{ Row = 24; D = "73.83";      E = "  +1.02%  " }
    @{ Row = 25; E = "  +0.15%  " }
    @{ Row = 26; D = "0.0000125";  E = "  +3.45%  " }
    @{ Row = 27; D = "9.57";       E = "  +0.84%  " }
    @{ Row = 28; D = "0.180";      E = "  -0.31%  " }
    @{ Row = 29; D = "0.998";      E = "  -0.19%  " }
    @{ Row = 30; D = "6.46";       E = "  +8.01%  " }
    @{ Row = 31; E = "  +4.59%  " }
    @{ Row = 32; D = "2.06";       E = "  +0.64%  " }
    @{ Row = 33; E = "  -1.12%  " }
    @{ Row = 34; D = "23.77";      E = "  -0.35%  " }
    @{ Row = 36; E = "  +2.18%  " }
    @{ Row = 37; E = "  +3.51%  " }
    @{ Row = 38; D = "161.84";     E = "  +0.38%  " }
    @{ Row = 39; D = "1.98";       E = "  +5.69%  " }
    @{ Row = 40; D = "3.052.65";   E = "  +5.31%  " }
    @{ Row = 41; D = "0.0775";     E = "  -1.10%  " }
    @{ Row = 42; D = "27.32";      E = "  -1.40%  " }
    @{ Row = 43; D = "0.0322";     E = "  -0.37%  " }
    @{ Row = 44; D = "4.53";       E = "  +2.61%  " }
    @{ Row = 45; D = "42.93";      E = "  +2.57%  " }
    @{ Row = 46; D = "0.778";      E = "  +0.11%  " }
    @{ Row = 47; D = "26.08";      E = "  +11.12%  " }
    @{ Row = 48; E = "  +2.71%  " }
    @{ Row = 49; D = "2.25";       E = "  +3.41%  " }
    @{ Row = 50; D = "6.73";       E = "  +2.65%  " }
    @{ Row = 51; D = "311.44";     E = "  +4.42%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey("D")) {
        $cell = $ws.Cells.Item($row, 4)
        # The "Price" column holds numeric-looking text (e.g. "580.50",
        # "65.734.19") that must stay literal text, matching the
        # original file's digits (incl. trailing zeros / dotted groups)
        # instead of being auto-coerced into a float by Excel's input
        # parser. Force text entry, then restore the original "Normal"
        # (general) cell style so no stray number-format is left behind.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }

    if ($u.ContainsKey("E")) {
        # "Volume(1h)" cells are already non-numeric strings
        # (padded with spaces and a %), so a plain assignment is safe.
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
